# Insert a new record row at row 448 (shifting the existing rows 448:485 down to 449:486)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(448).Insert()

$ws.Cells.Item(448, 1).Value  = 8
$ws.Cells.Item(448, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(448, 3).Value  = "Coquimbo"
$ws.Cells.Item(448, 4).Value  = 45265
$ws.Cells.Item(448, 5).Value  = 4
$ws.Cells.Item(448, 6).Value  = 100112031
$ws.Cells.Item(448, 7).Value  = "Poroto verde"
$ws.Cells.Item(448, 8).Value  = "Sin especificar"
$ws.Cells.Item(448, 9).Value  = "Primera"
$ws.Cells.Item(448, 10).Value = 520
$ws.Cells.Item(448, 11).Value = 31000
$ws.Cells.Item(448, 12).Value = 32000
$ws.Cells.Item(448, 13).Value = 31500
$ws.Cells.Item(448, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(448, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(448, 16).Value = 1260
$ws.Cells.Item(448, 17).Value = 25
$ws.Cells.Item(448, 18).Value = "Hortaliza"
